# Update "想去人数" (F) and "最低票价" (G) figures on the 展览 and 全部类型 sheets
# to reflect the latest scrape, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 286
$wsExpo.Range("F3").Value = 1169
$wsExpo.Range("F4").Value = 16652
$wsExpo.Range("F5").Value = 21
$wsExpo.Range("F6").Value = 1632
$wsExpo.Range("F7").Value = 60
$wsExpo.Range("F8").Value = 1
$wsExpo.Range("G8").Value = 39.9
$wsExpo.Range("F10").Value = 210
$wsExpo.Range("F12").Value = 11577
$wsExpo.Range("F14").Value = 1257
$wsExpo.Range("F15").Value = 4580
$wsExpo.Range("F16").Value = 416
$wsExpo.Range("F19").Value = 877
$wsExpo.Range("F20").Value = 333

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 286
$wsAll.Range("F4").Value = 1169
$wsAll.Range("F5").Value = 16652
$wsAll.Range("F6").Value = 21
$wsAll.Range("F7").Value = 1632
$wsAll.Range("F8").Value = 60
$wsAll.Range("F9").Value = 1
$wsAll.Range("G9").Value = 39.9
$wsAll.Range("F11").Value = 210
$wsAll.Range("F15").Value = 11577
$wsAll.Range("F17").Value = 1257
$wsAll.Range("F18").Value = 4580
$wsAll.Range("F19").Value = 416
$wsAll.Range("F22").Value = 877
$wsAll.Range("F23").Value = 333
